$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C2").Value = 12
$ws.Range("B3").Value = 4.5
$ws.Range("C3").Value = 12
$ws.Range("B4").Value = 0.75
$ws.Range("C5").Value = 20

# Update the selected cell (D6 -> D5)
$ws.Range("D5").Select()
